# Add 2022-Q3 data
# -----------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q3" positioned right after the
#    "总计" (summary) sheet and before the "2022-Q2" sheet.
# 2) Populate it with the quarterly fund-holding breakdown.
# 3) Insert a new row into the "总计" sheet for the 2022-Q3 totals,
#    shifting the existing rows down.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet    = $wb.Worksheets.Item(2)

# ---- 1) Create the new "2022-Q3" worksheet just before "2022-Q2" ----
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# ---- 2) Populate header row (B1:H1) ----
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3Sheet.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}

# ---- Fund rows (A2:H7) ----
# Columns: A index, B code, C name, D scale, E stock position, F position ratio,
#          G holding market value (100M yuan), H position rank
$rows = @(
    @(0, "006551", "中庚价值领航混合",             "118.19", "91.86", "3.39", "4.0066", 9),
    @(1, "011174", "中庚价值品质一年持有期混合",   "66.33",  "92.24", "2.57", "1.7047", 10),
    @(2, "007497", "中庚价值灵动灵活配置混合",     "36.46",  "89.30", "4.13", "1.5058", 3),
    @(3, "004702", "南方金融主题灵活配置混合A",    "11.74",  "92.15", "3.42", "0.4015", 8),
    @(4, "013500", "南方金融主题灵活配置混合C",    "7.39",   "92.15", "3.42", "0.2527", 8),
    @(5, "257040", "国联安红利混合",               "0.65",   "76.59", "5.40", "0.0351", 5)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $cellA = $q3Sheet.Range("A" + $r)
    $cellA.Value = $row[0]
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108   # xlCenter
    $cellA.VerticalAlignment = -4160     # xlTop
    $cellA.Borders.LineStyle = 1
    $cellA.Borders.Weight = 2

    $cellB = $q3Sheet.Range("B" + $r)
    $cellB.NumberFormat = "@"
    $cellB.Value = $row[1]

    $cellC = $q3Sheet.Range("C" + $r)
    $cellC.NumberFormat = "@"
    $cellC.Value = $row[2]

    $cellD = $q3Sheet.Range("D" + $r)
    $cellD.NumberFormat = "@"
    $cellD.Value = $row[3]

    $cellE = $q3Sheet.Range("E" + $r)
    $cellE.NumberFormat = "@"
    $cellE.Value = $row[4]

    $cellF = $q3Sheet.Range("F" + $r)
    $cellF.NumberFormat = "@"
    $cellF.Value = $row[5]

    $cellG = $q3Sheet.Range("G" + $r)
    $cellG.NumberFormat = "@"
    $cellG.Value = $row[6]

    $cellH = $q3Sheet.Range("H" + $r)
    $cellH.Value = $row[7]
}

$q3Sheet.Range("A1").Select()

# ---- 3) Update the "总计" (summary) sheet ----
# Insert a new row 2, shifting the existing data rows down by one.
$totalSheet.Range("A2:D2").Insert(-4121)   # xlShiftDown

# The inserted cells can inherit stray formatting from neighbouring rows;
# strip it so the new row 2 (B2:D2) starts out unformatted, matching the
# plain (non-bold, borderless) look of the other data cells.
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A2").Font.Bold = $true
$totalSheet.Range("A2").HorizontalAlignment = -4108
$totalSheet.Range("A2").VerticalAlignment = -4160
$totalSheet.Range("A2").Borders.LineStyle = 1
$totalSheet.Range("A2").Borders.Weight = 2

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 7.91

# Renumber the "A" index column (0-based sequential index) for the rows
# that were pushed down, and restore their values explicitly so we don't
# rely on whatever the Insert operation left behind.
$totalData = @(
    @(1, "2022-Q2", 15, 9.85),
    @(2, "2022-Q1", 24, 8.55),
    @(3, "2021-Q4", 9, 8.44),
    @(4, "2021-Q3", 8, 6.96),
    @(5, "2021-Q2", 3, 0.34)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $r = $i + 3
    $row = $totalData[$i]

    $cellA = $totalSheet.Range("A" + $r)
    $cellA.Value = $row[0]
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1
    $cellA.Borders.Weight = 2

    $totalSheet.Range("B" + $r).Value = $row[1]
    $totalSheet.Range("C" + $r).Value = $row[2]
    $totalSheet.Range("D" + $r).Value = $row[3]
}

$totalSheet.Range("A1").Select()
